$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.092652440071106
$ws.Range("B1").Value = 1.733769774436951
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.195434808731079
